$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($RowNum, $Values)
    for ($i = 0; $i -lt $Values.Count; $i++) {
        $ws.Cells.Item($RowNum, $i + 1).Value = $Values[$i]
    }
}

# Update the "data as of" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Agosto de 2020 a las 10:00"

# Row layout: Country, TotalCases, NewCases, ActiveCases, Recovered, Critical, NewDeaths, Deaths

# India (row 6) - values refreshed
Set-Row 6 @("India", 3170942, 6061, 2405047, 707325, 0, 24, 58570)

# Rusia (row 7) - values refreshed
Set-Row 7 @("Rusia", 966189, 4696, 779747, 169874, 0, 120, 16568)

# Ucrania overtakes Ecuador -> rows 30/31 swap order, new Ucrania data, Ecuador keeps its prior data
Set-Row 30 @("Ucrania", 108415, 1658, 52870, 53227, 0, 25, 2318)
Set-Row 31 @("Ecuador", 108289, 0, 94878, 7089, 0, 0, 6322)

# Singapur (row 49) - values refreshed
Set-Row 49 @("Singapur", 56435, 31, 54587, 1821, 0, 0, 27)

# Armenia (row 57) - values refreshed
Set-Row 57 @("Armenia", 42936, 111, 36475, 5603, 0, 4, 858)

# Maldivas (row 104) - values refreshed
Set-Row 104 @("Maldivas", 6912, 0, 4297, 2587, 0, 1, 28)

# Hungria (row 109) - values refreshed
Set-Row 109 @("Hungria", 5215, 24, 3716, 885, 0, 1, 614)

# Estonia overtakes Siria -> rows 134/135 swap order, new Estonia data, Siria keeps its prior data
Set-Row 134 @("Estonia", 2294, 19, 2038, 192, 0, 0, 64)
Set-Row 135 @("Siria", 2293, 0, 519, 1682, 0, 0, 92)

# Laos (row 206) - values refreshed
Set-Row 206 @("Laos", 22, 0, 21, 1, 0, 0, 0)

# Islas Malvinas overtakes Montserrat -> rows 214/215 swap order
Set-Row 214 @("Islas Malvinas", 13, 0, 13, 0, 0, 0, 0)
Set-Row 215 @("Montserrat", 13, 0, 12, 0, 0, 0, 1)
